$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.43773118810929
$ws.Range("B3").Value = 2.100966953615222
$ws.Range("B4").Value = 34.49896502285155
$ws.Range("B5").Value = 60
